# "Add files via upload" — add a new worksheet ("Sheet2") containing two more
# city names, place it after Sheet1, and make it the active/selected sheet
# (mirrors a user adding a second sheet of data and leaving it selected).

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no args inserts before the active sheet, so the new
# sheet starts out in front of Sheet1 -- move it after Sheet1 once created.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet2"

$newSheet.Range("A1").Value = "Kolkata"
$newSheet.Range("A2").Value = "Hyderabad"

# Put Sheet2 right after Sheet1 in tab order.
$newSheet.Move($null, $wb.Worksheets.Item("Sheet1"))

# Column A on the new sheet was sized to fit its (longer) city names.
$wb.Worksheets.Item("Sheet2").Columns("A:A").AutoFit()

# Leave the new sheet selected/active with B3 highlighted, matching the
# saved UI state.
$s2 = $wb.Worksheets.Item("Sheet2")
$s2.Range("B3").Select()
$s2.Activate()
